$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 672.05884
$ws.Range("I19").Value = 941.5
$ws.Range("J19").Value = 432.55554
$ws.Range("K19").Value = 941.5
$ws.Range("L19").Value = 432.55554
$ws.Range("M19").Value = -766.5
$ws.Range("N19").Value = -782.5555400000001
# Row 40
$ws.Range("H40").Value = 4180
$ws.Range("I40").Value = 4860
$ws.Range("J40").Value = 3160
$ws.Range("K40").Value = 4860
$ws.Range("L40").Value = 3160
$ws.Range("M40").Value = -4685
$ws.Range("N40").Value = -3510
# Row 51
$ws.Range("H51").Value = 4750
$ws.Range("J51").Value = 4300
$ws.Range("L51").Value = 4300
$ws.Range("N51").Value = -5268
# Row 64
$ws.Range("H64").Value = 3016.0625
$ws.Range("I64").Value = 2820.625
$ws.Range("K64").Value = 2820.625
$ws.Range("M64").Value = -2572.625
# Row 67
$ws.Range("H67").Value = 3016.0625
$ws.Range("I67").Value = 2820.625
$ws.Range("K67").Value = 2820.625
$ws.Range("M67").Value = -1962.625
# Row 112
$ws.Range("H112").Value = 1056.1786
$ws.Range("I112").Value = 766.1667
$ws.Range("J112").Value = 1135.2727
$ws.Range("K112").Value = 2298.5001
$ws.Range("L112").Value = 3405.8181
$ws.Range("M112").Value = -1190.5001
$ws.Range("N112").Value = -5621.8181
# Row 129
$ws.Range("H129").Value = 7464.273
$ws.Range("I129").Value = 472.5
$ws.Range("J129").Value = 9701.639999999999
$ws.Range("K129").Value = 1417.5
$ws.Range("L129").Value = 29104.92
$ws.Range("M129").Value = 3582.5
$ws.Range("N129").Value = -39104.92
# Row 137
$ws.Range("H137").Value = 1520.9656
$ws.Range("I137").Value = 1411.2778
$ws.Range("K137").Value = 4233.8334
$ws.Range("M137").Value = -1683.8334
# Row 138
$ws.Range("H138").Value = 2950.162
$ws.Range("I138").Value = 2035.25
$ws.Range("J138").Value = 3202.5518
$ws.Range("K138").Value = 6105.75
$ws.Range("L138").Value = 9607.6554
$ws.Range("M138").Value = -965.75
$ws.Range("N138").Value = -19887.6554

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 14
$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 2000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 2000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -2350
# Row 21
$ws.Range("H21").Value = 1372.3334
$ws.Range("I21").Value = 100
$ws.Range("J21").Value = 2008.5
$ws.Range("K21").Value = 100
$ws.Range("L21").Value = 2008.5
$ws.Range("M21").Value = 274
$ws.Range("N21").Value = -2756.5
# Row 29
$ws.Range("H29").Value = 2000
$ws.Range("I29").Value = 2000
$ws.Range("K29").Value = 2000
$ws.Range("M29").Value = -1692
# Row 30
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("K30").Value = 1000
$ws.Range("M30").Value = -850
# Row 32
$ws.Range("H32").Value = 15457.212
$ws.Range("I32").Value = 4291.129
$ws.Range("J32").Value = 34167.945
$ws.Range("K32").Value = 4291.129
$ws.Range("L32").Value = 34167.945
$ws.Range("M32").Value = -4004.129
$ws.Range("N32").Value = -34741.945
# Row 102
$ws.Range("H102").Value = 1432.2222
$ws.Range("I102").Value = 1398.8235
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1398.8235
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 223.1765
$ws.Range("N102").Value = -5244
# Row 122
$ws.Range("H122").Value = 1699.5769
$ws.Range("I122").Value = 916.7222
$ws.Range("J122").Value = 3461
$ws.Range("K122").Value = 2750.1666
$ws.Range("L122").Value = 10383
$ws.Range("M122").Value = -300.1666
$ws.Range("N122").Value = -15283

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 20002800
$ws.Range("I86").Value = 25002624
$ws.Range("J86").Value = 3495
$ws.Range("K86").Value = 25002624
$ws.Range("L86").Value = 3495
$ws.Range("M86").Value = -25001501
$ws.Range("N86").Value = -5741
# Row 89
$ws.Range("H89").Value = 20002800
$ws.Range("I89").Value = 25002624
$ws.Range("J89").Value = 3495
$ws.Range("K89").Value = 125013120
$ws.Range("L89").Value = 17475
$ws.Range("M89").Value = -125007504
$ws.Range("N89").Value = -28707

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 83336120
$ws.Range("J107").Value = 2204.3333
$ws.Range("L107").Value = 2204.3333
$ws.Range("N107").Value = -6044.3333

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 726.125
$ws.Range("J122").Value = 1198.8667
$ws.Range("L122").Value = 10789.8003
$ws.Range("N122").Value = -15689.8003

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 634.6
$ws.Range("I9").Value = 634.6
$ws.Range("K9").Value = 634.6
$ws.Range("M9").Value = -464.6
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
# Row 113
$ws.Range("H113").Value = 1475.7
$ws.Range("I113").Value = 1091.5
$ws.Range("K113").Value = 1091.5
$ws.Range("M113").Value = 1078.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 770.5
$ws.Range("I16").Value = 818.13336
$ws.Range("J16").Value = 627.6
$ws.Range("K16").Value = 818.13336
$ws.Range("L16").Value = 627.6
$ws.Range("M16").Value = -648.13336
$ws.Range("N16").Value = -967.6
# Row 34
$ws.Range("H34").Value = 12266.667
$ws.Range("I34").Value = 10400
$ws.Range("J34").Value = 16000
$ws.Range("K34").Value = 10400
$ws.Range("L34").Value = 16000
$ws.Range("M34").Value = -10228
$ws.Range("N34").Value = -16344
# Row 46
$ws.Range("H46").Value = 971.6429000000001
$ws.Range("I46").Value = 838.5
$ws.Range("J46").Value = 1024.9
$ws.Range("K46").Value = 838.5
$ws.Range("L46").Value = 1024.9
$ws.Range("M46").Value = -650.5
$ws.Range("N46").Value = -1400.9
# Row 61
$ws.Range("H61").Value = 2101
$ws.Range("I61").Value = 1692.25
$ws.Range("K61").Value = 1692.25
$ws.Range("M61").Value = -1490.25
# Row 113
$ws.Range("H113").Value = 2101
$ws.Range("I113").Value = 1692.25
$ws.Range("K113").Value = 1692.25
$ws.Range("M113").Value = 477.75

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 32
$ws.Range("H32").Value = 1800
$ws.Range("I32").Value = 1800
$ws.Range("K32").Value = 1800
$ws.Range("M32").Value = -1483
